# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-suffix columns to "_FV2210" / "_FV2304"
# - Turn the data range into a proper Excel Table (Table1)
# - Freeze the header row (split/freeze pane below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used      = $ws.UsedRange
$lastRow   = $used.Rows.Count
$lastCol   = $used.Columns.Count

# 1. Rename header row cells: *_old -> *_FV2210, *_new -> *_FV2304
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $name = $cell.Value()
    if ($name -ne $null) {
        if ($name.EndsWith("_old")) {
            $cell.Value = $name.Substring(0, $name.Length - 4) + "_FV2210"
        }
        elseif ($name.EndsWith("_new")) {
            $cell.Value = $name.Substring(0, $name.Length - 4) + "_FV2304"
        }
    }
}

# 2. Convert the data range into a proper Excel Table named "Table1"
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# 3. Freeze panes below the header row (select A2, then freeze)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
